# [Word] (TableCell) Map existing sample
# Adds a new "TableCell" row to the Snippets table (row 70), mirroring the
# existing "word-tables-table-cell-access" / "getTableCell" sample that's
# already referenced by the "Table" class (row 69), but for the "TableCell"
# class / "class" member.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the "Snippets" table by one row - this keeps the table ref,
# autoFilter ref and sheet dimension in sync automatically.
$tbl = $ws.ListObjects.Item("Snippets")
$newRow = $tbl.ListRows.Add()

$newRowIndex = $tbl.Range.Rows.Count
$rowNum = $tbl.Range.Row + $newRowIndex - 1

# Fill in the new row's data.
$ws.Cells.Item($rowNum, 1).Value = "TableCell"
$ws.Cells.Item($rowNum, 2).Value = ""
$ws.Cells.Item($rowNum, 3).Value = "class"
$ws.Cells.Item($rowNum, 4).Value = "word-tables-table-cell-access"
$ws.Cells.Item($rowNum, 5).Value = "getTableCell"

# Columns A:B on the new row use a plain "General" number format (clears any
# inherited formatting from the table), matching how the row was entered.
$ws.Range("A" + $rowNum + ":B" + $rowNum).NumberFormat = "General"

# Columns D:E pick up the same look as the row above (vertical-centered
# font) that the rest of the table body uses.
$ws.Range("D" + ($rowNum - 1) + ":E" + ($rowNum - 1)).Copy() | Out-Null
$ws.Range("D" + $rowNum + ":E" + $rowNum).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Move the selection to the newly-added cell, like a user would after typing
# the last value in the row.
$ws.Range("E" + $rowNum).Select()
